$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for rows 4-12 (columns D, I, J, K, L, M, P)
$rows = @{
    4  = @{ D = 44536; I = "Primera"; J = 87; K = 22000; L = 22000; M = 22000; P = 1222 }
    5  = @{ D = 44536; I = "Segunda"; J = 80; K = 20000; L = 20000; M = 20000; P = 1111 }
    6  = @{ D = 44242; I = "Primera"; J = 60; K = 13000; L = 13000; M = 13000; P = 722 }
    7  = @{ D = 44242; I = "Segunda"; J = 50; K = 10000; L = 10000; M = 10000; P = 556 }
    8  = @{ D = 44235; I = "Primera"; J = 80; K = 14000; L = 14000; M = 14000; P = 778 }
    9  = @{ D = 44235; I = "Segunda"; J = 70; K = 12000; L = 12000; M = 12000; P = 667 }
    10 = @{ D = 44238; I = "Tercera"; J = 60; K = 10000; L = 10000; M = 10000; P = 556 }
    11 = @{ D = 44238; I = "Primera"; J = 90; K = 13000; L = 13000; M = 13000; P = 722 }
    12 = @{ D = 44238; I = "Segunda"; J = 80; K = 11000; L = 11000; M = 11000; P = 611 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("I$r").Value = $vals.I
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
